# Auto-generated script to apply numeric corrections to Aegis_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 1276.72
$ws.Cells.Item(33, 10).Value = 4519.6
$ws.Cells.Item(33, 12).Value = 4519.6
$ws.Cells.Item(33, 14).Value = -4977.6
$ws.Cells.Item(42, 8).Value = 2083721.1
$ws.Cells.Item(42, 10).Value = 594.3333
$ws.Cells.Item(42, 12).Value = 1782.9999
$ws.Cells.Item(42, 14).Value = -2242.9999
$ws.Cells.Item(43, 8).Value = 1653.5927
$ws.Cells.Item(43, 10).Value = 1340.4166
$ws.Cells.Item(43, 12).Value = 1340.4166
$ws.Cells.Item(43, 14).Value = -1478.4166
$ws.Cells.Item(48, 8).Value = 2750
$ws.Cells.Item(48, 10).Value = 2857.1428
$ws.Cells.Item(48, 12).Value = 8571.428400000001
$ws.Cells.Item(48, 14).Value = -9155.428400000001
$ws.Cells.Item(56, 8).Value = 2750
$ws.Cells.Item(56, 10).Value = 2857.1428
$ws.Cells.Item(56, 12).Value = 8571.428400000001
$ws.Cells.Item(56, 14).Value = -9639.428400000001
$ws.Cells.Item(86, 8).Value = 3968.087
$ws.Cells.Item(86, 9).Value = 1411.25
$ws.Cells.Item(86, 10).Value = 5331.7334
$ws.Cells.Item(86, 11).Value = 1411.25
$ws.Cells.Item(86, 12).Value = 5331.7334
$ws.Cells.Item(86, 13).Value = -288.25
$ws.Cells.Item(86, 14).Value = -7577.7334
$ws.Cells.Item(89, 8).Value = 3968.087
$ws.Cells.Item(89, 9).Value = 1411.25
$ws.Cells.Item(89, 10).Value = 5331.7334
$ws.Cells.Item(89, 11).Value = 7056.25
$ws.Cells.Item(89, 12).Value = 26658.667
$ws.Cells.Item(89, 13).Value = -1440.25
$ws.Cells.Item(89, 14).Value = -37890.667
$ws.Cells.Item(125, 8).Value = 2972.318
$ws.Cells.Item(125, 9).Value = 2617.8462
$ws.Cells.Item(125, 10).Value = 3484.3333
$ws.Cells.Item(125, 11).Value = 23560.6158
$ws.Cells.Item(125, 12).Value = 31358.9997
$ws.Cells.Item(125, 13).Value = -21100.6158
$ws.Cells.Item(125, 14).Value = -36278.9997
$ws.Cells.Item(129, 8).Value = 4885.346
$ws.Cells.Item(129, 10).Value = 1236.1
$ws.Cells.Item(129, 12).Value = 3708.3
$ws.Cells.Item(129, 14).Value = -13708.3
$ws.Cells.Item(131, 8).Value = 4614.5
$ws.Cells.Item(131, 10).Value = 4944.9736
$ws.Cells.Item(131, 12).Value = 14834.9208
$ws.Cells.Item(131, 14).Value = -24914.9208

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 25265.791
$ws.Cells.Item(32, 9).Value = 4331.4067
$ws.Cells.Item(32, 10).Value = 120275.695
$ws.Cells.Item(32, 11).Value = 4331.4067
$ws.Cells.Item(32, 12).Value = 120275.695
$ws.Cells.Item(32, 13).Value = -4044.4067
$ws.Cells.Item(32, 14).Value = -120849.695
$ws.Cells.Item(45, 8).Value = 500406
$ws.Cells.Item(45, 9).Value = 500406
$ws.Cells.Item(45, 11).Value = 500406
$ws.Cells.Item(45, 13).Value = -500029
$ws.Cells.Item(61, 8).Value = 1970.2826
$ws.Cells.Item(61, 9).Value = 1431.5
$ws.Cells.Item(61, 10).Value = 2670.7
$ws.Cells.Item(61, 11).Value = 1431.5
$ws.Cells.Item(61, 12).Value = 2670.7
$ws.Cells.Item(61, 13).Value = -1219.5
$ws.Cells.Item(61, 14).Value = -3094.7
$ws.Cells.Item(97, 8).Value = 85133.25
$ws.Cells.Item(97, 9).Value = 91862.55
$ws.Cells.Item(97, 11).Value = 91862.55
$ws.Cells.Item(97, 13).Value = -91366.55
$ws.Cells.Item(122, 8).Value = 1130.5667
$ws.Cells.Item(122, 9).Value = 1024.5385
$ws.Cells.Item(122, 10).Value = 1819.75
$ws.Cells.Item(122, 11).Value = 3073.6155
$ws.Cells.Item(122, 12).Value = 5459.25
$ws.Cells.Item(122, 13).Value = -623.6155000000003
$ws.Cells.Item(122, 14).Value = -10359.25
$ws.Cells.Item(132, 8).Value = 2462.5898
$ws.Cells.Item(132, 9).Value = 2207.147
$ws.Cells.Item(132, 10).Value = 4199.6
$ws.Cells.Item(132, 11).Value = 6621.441
$ws.Cells.Item(132, 12).Value = 12598.8
$ws.Cells.Item(132, 13).Value = -4091.441
$ws.Cells.Item(132, 14).Value = -17658.8
$ws.Cells.Item(136, 8).Value = 1970.2826
$ws.Cells.Item(136, 9).Value = 1431.5
$ws.Cells.Item(136, 10).Value = 2670.7
$ws.Cells.Item(136, 11).Value = 4294.5
$ws.Cells.Item(136, 12).Value = 8012.099999999999
$ws.Cells.Item(136, 13).Value = -1744.5
$ws.Cells.Item(136, 14).Value = -13112.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 42919.375
$ws.Cells.Item(20, 9).Value = 51232.35
$ws.Cells.Item(20, 10).Value = 1354.5
$ws.Cells.Item(20, 11).Value = 51232.35
$ws.Cells.Item(20, 12).Value = 1354.5
$ws.Cells.Item(20, 13).Value = -50985.35
$ws.Cells.Item(20, 14).Value = -1848.5
$ws.Cells.Item(94, 8).Value = 668.6667
$ws.Cells.Item(94, 9).Value = 574
$ws.Cells.Item(94, 11).Value = 574
$ws.Cells.Item(94, 13).Value = -123
$ws.Cells.Item(107, 8).Value = 66721096
$ws.Cells.Item(107, 9).Value = 71486744
$ws.Cells.Item(107, 11).Value = 71486744
$ws.Cells.Item(107, 13).Value = -71484824
$ws.Cells.Item(134, 8).Value = 3208.9644
$ws.Cells.Item(134, 9).Value = 3272.2593
$ws.Cells.Item(134, 10).Value = 1500
$ws.Cells.Item(134, 11).Value = 9816.777900000001
$ws.Cells.Item(134, 12).Value = 4500
$ws.Cells.Item(134, 13).Value = -7281.777900000001
$ws.Cells.Item(134, 14).Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 33860.47
$ws.Cells.Item(31, 9).Value = 1178.2413
$ws.Cells.Item(31, 10).Value = 81249.7
$ws.Cells.Item(31, 11).Value = 1178.2413
$ws.Cells.Item(31, 12).Value = 81249.7
$ws.Cells.Item(31, 13).Value = -883.2412999999999
$ws.Cells.Item(31, 14).Value = -81839.7
$ws.Cells.Item(34, 8).Value = 33860.47
$ws.Cells.Item(34, 9).Value = 1178.2413
$ws.Cells.Item(34, 10).Value = 81249.7
$ws.Cells.Item(34, 11).Value = 1178.2413
$ws.Cells.Item(34, 12).Value = 81249.7
$ws.Cells.Item(34, 13).Value = -976.2412999999999
$ws.Cells.Item(34, 14).Value = -81653.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 191.75
$ws.Cells.Item(8, 9).Value = 191.75
$ws.Cells.Item(8, 11).Value = 575.25
$ws.Cells.Item(8, 13).Value = -436.25
$ws.Cells.Item(37, 8).Value = 24942.75
$ws.Cells.Item(37, 10).Value = 24942.75
$ws.Cells.Item(37, 12).Value = 74828.25
$ws.Cells.Item(37, 14).Value = -75052.25
$ws.Cells.Item(69, 8).Value = 2742.8
$ws.Cells.Item(69, 10).Value = 2928.5
$ws.Cells.Item(69, 12).Value = 8785.5
$ws.Cells.Item(69, 14).Value = -10407.5
$ws.Cells.Item(72, 8).Value = 2742.8
$ws.Cells.Item(72, 10).Value = 2928.5
$ws.Cells.Item(72, 12).Value = 26356.5
$ws.Cells.Item(72, 14).Value = -34468.5
$ws.Cells.Item(80, 8).Value = 12807.223
$ws.Cells.Item(80, 9).Value = 999
$ws.Cells.Item(80, 10).Value = 14283.25
$ws.Cells.Item(80, 11).Value = 2997
$ws.Cells.Item(80, 12).Value = 42849.75
$ws.Cells.Item(80, 13).Value = -2061
$ws.Cells.Item(80, 14).Value = -44721.75
$ws.Cells.Item(83, 8).Value = 12807.223
$ws.Cells.Item(83, 9).Value = 999
$ws.Cells.Item(83, 10).Value = 14283.25
$ws.Cells.Item(83, 11).Value = 8991
$ws.Cells.Item(83, 12).Value = 128549.25
$ws.Cells.Item(83, 13).Value = -4311
$ws.Cells.Item(83, 14).Value = -137909.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(23, 8).Value = 18428.715
$ws.Cells.Item(23, 10).Value = 25799.8
$ws.Cells.Item(23, 12).Value = 25799.8
$ws.Cells.Item(23, 14).Value = -26245.8
$ws.Cells.Item(80, 8).Value = 100101390
$ws.Cells.Item(80, 9).Value = 250250670
$ws.Cells.Item(80, 10).Value = 1864.1666
$ws.Cells.Item(80, 11).Value = 250250670
$ws.Cells.Item(80, 12).Value = 1864.1666
$ws.Cells.Item(80, 13).Value = -250249672
$ws.Cells.Item(80, 14).Value = -3860.1666
$ws.Cells.Item(83, 8).Value = 100101390
$ws.Cells.Item(83, 9).Value = 250250670
$ws.Cells.Item(83, 10).Value = 1864.1666
$ws.Cells.Item(83, 11).Value = 1251253350
$ws.Cells.Item(83, 12).Value = 9320.833000000001
$ws.Cells.Item(83, 13).Value = -1251248358
$ws.Cells.Item(83, 14).Value = -19304.833
$ws.Cells.Item(97, 8).Value = 166670780
$ws.Cells.Item(97, 9).Value = 166670780
$ws.Cells.Item(97, 11).Value = 166670780
$ws.Cells.Item(97, 13).Value = -166670284
$ws.Cells.Item(120, 8).Value = 34758.25
$ws.Cells.Item(120, 10).Value = 34758.25
$ws.Cells.Item(120, 12).Value = 34758.25
$ws.Cells.Item(120, 14).Value = -44434.25
$ws.Cells.Item(122, 8).Value = 2034.1333
$ws.Cells.Item(122, 9).Value = 2036.5714
$ws.Cells.Item(122, 11).Value = 6109.7142
$ws.Cells.Item(122, 13).Value = -3659.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1049.0555
$ws.Cells.Item(22, 10).Value = 1034.5714
$ws.Cells.Item(22, 12).Value = 1034.5714
$ws.Cells.Item(22, 14).Value = -1624.5714
$ws.Cells.Item(26, 8).Value = 15000
$ws.Cells.Item(26, 9).Value = 10000
$ws.Cells.Item(26, 10).Value = 20000
$ws.Cells.Item(26, 11).Value = 10000
$ws.Cells.Item(26, 12).Value = 20000
$ws.Cells.Item(26, 13).Value = -9705
$ws.Cells.Item(26, 14).Value = -20590
$ws.Cells.Item(27, 8).Value = 1049.0555
$ws.Cells.Item(27, 10).Value = 1034.5714
$ws.Cells.Item(27, 12).Value = 1034.5714
$ws.Cells.Item(27, 14).Value = -1248.5714
$ws.Cells.Item(76, 8).Value = 288
$ws.Cells.Item(76, 10).Value = 288
$ws.Cells.Item(76, 12).Value = 288
$ws.Cells.Item(76, 14).Value = -964
$ws.Cells.Item(79, 8).Value = 288
$ws.Cells.Item(79, 10).Value = 288
$ws.Cells.Item(79, 12).Value = 288
$ws.Cells.Item(79, 14).Value = -2628
$ws.Cells.Item(122, 8).Value = 2071
$ws.Cells.Item(122, 9).Value = 2004
$ws.Cells.Item(122, 11).Value = 6012
$ws.Cells.Item(122, 13).Value = -3562
$ws.Cells.Item(132, 8).Value = 3578.8333
$ws.Cells.Item(132, 9).Value = 3947.8076
$ws.Cells.Item(132, 11).Value = 11843.4228
$ws.Cells.Item(132, 13).Value = -9313.4228
$ws.Cells.Item(136, 8).Value = 2099.1667
$ws.Cells.Item(136, 9).Value = 2076.3635
$ws.Cells.Item(136, 10).Value = 2350
$ws.Cells.Item(136, 11).Value = 6229.0905
$ws.Cells.Item(136, 12).Value = 7050
$ws.Cells.Item(136, 13).Value = -3679.0905
$ws.Cells.Item(136, 14).Value = -12150

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 286699.84
$ws.Cells.Item(81, 9).Value = 250875
$ws.Cells.Item(81, 10).Value = 334466.34
$ws.Cells.Item(81, 11).Value = 501750
$ws.Cells.Item(81, 12).Value = 668932.6800000001
$ws.Cells.Item(81, 13).Value = -500689
$ws.Cells.Item(81, 14).Value = -671054.6800000001
$ws.Cells.Item(84, 8).Value = 286699.84
$ws.Cells.Item(84, 9).Value = 250875
$ws.Cells.Item(84, 10).Value = 334466.34
$ws.Cells.Item(84, 11).Value = 2508750
$ws.Cells.Item(84, 12).Value = 3344663.4
$ws.Cells.Item(84, 13).Value = -2503446
$ws.Cells.Item(84, 14).Value = -3355271.4
$ws.Cells.Item(126, 8).Value = 1668.4667
$ws.Cells.Item(126, 9).Value = 1522.4546
$ws.Cells.Item(126, 10).Value = 2070
$ws.Cells.Item(126, 11).Value = 4567.3638
$ws.Cells.Item(126, 12).Value = 6210
$ws.Cells.Item(126, 13).Value = -2097.3638
$ws.Cells.Item(126, 14).Value = -11150
$ws.Cells.Item(132, 8).Value = 5050.421
$ws.Cells.Item(132, 9).Value = 6838.125
$ws.Cells.Item(132, 10).Value = 3750.2727
$ws.Cells.Item(132, 11).Value = 20514.375
$ws.Cells.Item(132, 12).Value = 11250.8181
$ws.Cells.Item(132, 13).Value = -17984.375
$ws.Cells.Item(132, 14).Value = -16310.8181
$ws.Cells.Item(135, 8).Value = 34140.273
$ws.Cells.Item(135, 10).Value = 34140.273
$ws.Cells.Item(135, 12).Value = 34140.273
$ws.Cells.Item(135, 14).Value = -44280.273
$ws.Cells.Item(140, 8).Value = 64857.25
$ws.Cells.Item(140, 10).Value = 64857.25
$ws.Cells.Item(140, 12).Value = 64857.25
$ws.Cells.Item(140, 14).Value = -75217.25
